$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 19.60879629629628
$ws.Range("R2").Value = 1.523070892784755
$ws.Range("S2").Value = 1.622381208625714

$ws.Range("K3").Value = 19.60879629629628

$ws.Range("K6").Value = 17.72453703703706

$ws.Range("K7").Value = -0.763888888888889
$ws.Range("R7").Value = 1.32738293362297
$ws.Range("S7").Value = 1.399902772843707

$ws.Range("K8").Value = -0.763888888888889

$ws.Range("K9").Value = 15.74228395061728

$ws.Range("K10").Value = 15.74228395061728
$ws.Range("R10").Value = 0.9359630391213685
$ws.Range("S10").Value = 0.9581010400034263

$ws.Range("K11").Value = -3.083333333333334

$ws.Range("K12").Value = 1.925925925925943

$ws.Range("K13").Value = 1.925925925925943

$ws.Range("K14").Value = 14.96875

$ws.Range("K15").Value = 14.96875
$ws.Range("R15").Value = 0.9337049960039765
$ws.Range("S15").Value = 0.9557151767570068

$ws.Range("K22").Value = 15.74228395061728

$ws.Range("K23").Value = 15.74228395061728
$ws.Range("R23").Value = 0.9359630391213685
$ws.Range("S23").Value = 0.9581010400034263

$ws.Range("K28").Value = 14.96875
$ws.Range("R28").Value = 1.473592088566053
$ws.Range("S28").Value = 1.565708370582976

$ws.Range("K29").Value = 14.96875

$ws.Range("K30").Value = 1.925925925925943

$ws.Range("K37").Value = 14.96875
$ws.Range("R37").Value = 0.9337049960039765
$ws.Range("S37").Value = 0.9557151767570068

$ws.Range("K38").Value = 14.96875
